$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column widths for E and G (closest values the engine's
# character-width grid can represent to the authored widths of
# 44.44140625 / 57.77734375)
$ws.Columns.Item(5).ColumnWidth = 43.6
$ws.Columns.Item(7).ColumnWidth = 57

# Week header for column G (week 44)
$ws.Range("G1").Value = 44

# New task assignments for week 44 (column G)
# Set in this order so new shared-string entries are appended in the
# same sequence as the target workbook (index 12,13,14,15).
$ws.Range("G5").Value = " State Diagram for Pedestrian And Car TLS And Bus Button"
$ws.Range("G7").Value = "Class Diagram for Pedestrian And Car TLS And Bus Button"
$ws.Range("G4").Value = "Coding for Arduino Circuit Pedestrian &Car traffic light And Bus Button"
$ws.Range("G6").Value = "Coding and connecting circuit through TinkerCAD And Bus Button"

# Update selection
$ws.Range("G13").Select()
